$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-34 (inclusive): column C -> 7811
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 3).Value = 7811
}

# Rows 35-252 (inclusive): column C -> 7622
for ($r = 35; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7622
}
